$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "link" column (P) previously held a "view" URL for every dataset, duplicating
# the real data-path link already present in column Q ("data_path"). Per the commit
# "fix: dataset link only for dataset with preview", only the dataset that actually
# has a preview (row 3 / accident_route) should keep a link, and it now points to the
# local preview file instead of an external "view" URL. All other rows get cleared.

$ws.Range("P2").ClearContents()
$ws.Range("P3").Value = "dataset/accident_route.xlsx"
$ws.Range("P4").ClearContents()
$ws.Range("P5").ClearContents()
$ws.Range("P6").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("P10").ClearContents()
$ws.Range("P11").ClearContents()
$ws.Range("P12").ClearContents()
$ws.Range("P13").ClearContents()
$ws.Range("P14").ClearContents()
$ws.Range("P15").ClearContents()
$ws.Range("P17").ClearContents()
$ws.Range("P18").ClearContents()
$ws.Range("P19").ClearContents()
$ws.Range("P20").ClearContents()
$ws.Range("P21").ClearContents()

